$wb = $excel.ActiveWorkbook

# --- Sheet "基金受益憑證" (fund / trust-certificate holdings, sheet4.xml) ---
# Row 1 currently just repeats the data values instead of holding proper
# column headers. Bring it in line with the header layout used by the other
# property sheets (name/owner/dealer/quantity/face_value/currency/total),
# and extend row 2 with the common property_category..index metadata
# columns that every other sheet already carries.
$ws4 = $wb.Worksheets.Item("基金受益憑證")

$ws4.Range("B1").Value = "name"
$ws4.Range("C1").Value = "owner"
$ws4.Range("D1").Value = "dealer"
$ws4.Range("E1").Value = "quantity"
$ws4.Range("F1").Value = "face_value"
$ws4.Range("G1").Value = "currency"
$ws4.Range("H1").Value = "total"
$ws4.Range("I1").Value = "property_category"
$ws4.Range("J1").Value = "category"
$ws4.Range("K1").Value = "date"
$ws4.Range("L1").Value = "legislator_name"
$ws4.Range("M1").Value = "legislator_id"
$ws4.Range("N1").Value = "source_file"
$ws4.Range("O1").Value = "index"

$ws4.Range("I2").Value = "fund"
$ws4.Range("J2").Value = "normal"

# "2012-04-30" must land as literal text (matching the rest of the
# workbook), but assigning it straight to .Value lets Excel's date
# auto-recognition turn it into a serial date. Route it through a text
# formula first, then collapse the formula down to its cached value with
# a values-only paste so the cell ends up a plain shared string again.
$ws4.Range("K2").Formula = "=""2012-04-30"""
$ws4.Range("K2").Copy()
$ws4.Range("K2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false

$ws4.Range("L2").Value = "洪秀柱"
$ws4.Range("M2").Value = 546
$ws4.Range("N2").Value = "tmp31791"
$ws4.Range("O2").Value = 67

# Match the header/data cell formatting (bold+bordered header row, plain
# data row) already used by B1:H1 / B2:H2 onto the newly added columns.
$ws4.Range("B1").Copy()
$ws4.Range("I1:O1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws4.Range("B2").Copy()
$ws4.Range("I2:O2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

Write-Output "done"
